$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = 0
$ws.Range("F7").Value = 2
$ws.Range("F9").Value = 5
$ws.Range("F13").Value = 2
$ws.Range("F16").Value = -2
$ws.Range("F17").Value = -1
$ws.Range("F20").Value = -2
$ws.Range("F28").Value = -4
